# Applies the "Added new sheet : LoginCredentials" edit to UserDetails.xlsx
# (credentials.xlsx fixture): fixes a couple of data-entry glitches on the
# first sheet, adds two new rows of users, renames the first sheet, and
# adds a second "LoginCredentials" sheet that mirrors the e-mail/password
# pairs from the first sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet1 ("credentials" -> "UserInformation") data fixes
# ---------------------------------------------------------------------

# Row 3 was "Subha" with no initials column value -> becomes "Ultimate"/"subha"
$ws1.Range("A3").Value = "Ultimate"
$ws1.Range("B3").Value = "subha"

# Row 2 gains an initials value in column B (was a stray numeric 21)
$ws1.Range("B2").Value = "pal"
$ws1.Range("B2").HorizontalAlignment = -4131   # xlHAlignLeft
$ws1.Range("B2").VerticalAlignment = -4160     # xlVAlignTop

$ws1.Range("B3").HorizontalAlignment = -4131   # xlHAlignLeft

# Row4's stray numeric 21 is cleared out entirely (kept left-aligned/empty)
$ws1.Range("B4").ClearContents()
$ws1.Range("B4").HorizontalAlignment = -4131   # xlHAlignLeft

# Row5's stray numeric 21 is removed completely (no style left behind)
$ws1.Range("B5").ClearContents()

# Row6 gains initials "tk"
$ws1.Range("B6").Value = "tk"

# New row 7: Rj / mohan / rj@bksoft.com / spdf@*9 / Yes
$ws1.Range("A7").Value = "Rj"
$ws1.Range("B7").Value = "mohan"
$ws1.Range("C7").Value = "rj@bksoft.com"
$ws1.Range("D7").Value = "spdf@*9"
$ws1.Range("E7").Value = "Yes"
$ws1.Hyperlinks.Add($ws1.Range("C7"), "mailto:rj@bksoft.com") | Out-Null
$ws1.Range("C7").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("D7"), "mailto:spdf@*9") | Out-Null
$ws1.Range("D7").Style = "Hyperlink"

# New row 8: Sb / tk / sb@bksoft.com / sb@bk / Yes
$ws1.Range("A8").Value = "Sb"
$ws1.Range("B8").Value = "tk"
$ws1.Range("C8").Value = "sb@bksoft.com"
$ws1.Range("D8").Value = "sb@bk"
$ws1.Range("E8").Value = "Yes"
$ws1.Hyperlinks.Add($ws1.Range("C8"), "mailto:sb@bksoft.com") | Out-Null
$ws1.Range("C8").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("D8"), "mailto:sb@bk") | Out-Null
$ws1.Range("D8").Style = "Hyperlink"

# Selection left on the e-mail/password block
$ws1.Range("C2:D8").Select()

# ---------------------------------------------------------------------
# New sheet: LoginCredentials (UserName / Password lookup table)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LoginCredentials"

$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"

$ws2.Range("A2").Value = "abhi@gmail.com"
$ws2.Range("B2").Value = "pal@123"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:abhi@gmail.com") | Out-Null
$ws2.Range("A2").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:pal@123") | Out-Null
$ws2.Range("B2").Style = "Hyperlink"

$ws2.Range("A3").Value = "subha@yahoo.com"
$ws2.Range("B3").Value = "sangita&45"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:subha@yahoo.com") | Out-Null
$ws2.Range("A3").Style = "Hyperlink"

$ws2.Range("A4").Value = "shiv@radiffmail.com"
$ws2.Range("B4").Value = "shiv%88"
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:shiv@radiffmail.com") | Out-Null
$ws2.Range("A4").Style = "Hyperlink"

$ws2.Range("A5").Value = "biky@gmail.com"
$ws2.Range("B5").Value = "kundu#56"
$ws2.Hyperlinks.Add($ws2.Range("A5"), "mailto:biky@gmail.com") | Out-Null
$ws2.Range("A5").Style = "Hyperlink"

$ws2.Range("A6").Value = "pupu@gmail.com"
$ws2.Range("B6").Value = "pup&90"
$ws2.Hyperlinks.Add($ws2.Range("A6"), "mailto:pupu@gmail.com") | Out-Null
$ws2.Range("A6").Style = "Hyperlink"

$ws2.Range("A7").Value = "rj@bksoft.com"
$ws2.Range("B7").Value = "spdf@*9"
$ws2.Hyperlinks.Add($ws2.Range("A7"), "mailto:rj@bksoft.com") | Out-Null
$ws2.Range("A7").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("B7"), "mailto:spdf@*9") | Out-Null
$ws2.Range("B7").Style = "Hyperlink"

$ws2.Range("A8").Value = "sb@bksoft.com"
$ws2.Range("B8").Value = "sb@bk"
$ws2.Hyperlinks.Add($ws2.Range("A8"), "mailto:sb@bksoft.com") | Out-Null
$ws2.Range("A8").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("B8"), "mailto:sb@bk") | Out-Null
$ws2.Range("B8").Style = "Hyperlink"

# Column widths to roughly match the "best fit" sizing of the source sheet
$ws2.Columns("A").ColumnWidth = 18.88
$ws2.Columns("B").ColumnWidth = 9.88

# Selection / activation: LoginCredentials becomes the active (front) tab
$ws2.Range("B4").Select()
$ws2.Activate()
$excel.ActiveWindow.Zoom = 175

# ---------------------------------------------------------------------
# Rename the first sheet last (keeps $ws1 easy to reference above)
# ---------------------------------------------------------------------
$ws1.Name = "UserInformation"

Write-Host "LoginCredentials sheet added"
